# device_route.xlsx: refresh the sample "line" device rows (H2:H4 / G2:G4)
# and set the selection/print setup to match the author's re-upload.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4 used to hold sample device names (hk158, crtech_1, crtech_2) with
# a stray "Nocnoc1123" username in column G. Replace them with the generic
# placeholder data (root / line1 / line2 / line3) that the new upload uses.
$ws.Range("G2").Value = "root"
$ws.Range("H2").Value = "line1"

$ws.Range("G3").Value = "root"
$ws.Range("H3").Value = "line2"

$ws.Range("G4").Value = "root"
$ws.Range("H4").Value = "line3"

# Selection now spans the edited H column cells.
$ws.Range("H2:H4").Select() | Out-Null

# Page setup: A4, portrait (adds a <pageSetup> element on save).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Best-effort: match the saved window geometry from the re-upload.
$excel.ActiveWindow.Left = 4695
$excel.ActiveWindow.Top = 4050
$excel.ActiveWindow.Width = 21600
$excel.ActiveWindow.Height = 12855
